$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeiterfassung")

# Row 11: meeting on 2024-10-22, 1h, "Aufgaben" / "Visionspowerpoint"
# Reuse the date format from the existing date column (A7:A10) via copy/paste of formats
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A11").Value = 45587
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "Aufgaben"
$ws.Range("D11").Value = "Visionspowerpoint"

# Row 12: meeting on 2024-10-23, 2h, "Aufgaben" / "Visionspowerpoint"
# This entry gets a new date format (d-mmm) - different from the previous rows
$ws.Range("A12").Value = 45588
$ws.Range("A12").NumberFormat = "d-mmm"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "Aufgaben"
$ws.Range("D12").Value = "Visionspowerpoint"

$excel.CutCopyMode = 0

$ws.Range("G11").Select() | Out-Null
